$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new shared string "Constant" at A3, and value 12 at B3
$ws.Range("A3").Value = "Constant"
$ws.Range("B3").Value = 12

# Update B2 formula to include reference to B3
$ws.Range("B2").Formula = "=EXP(-1/input)+(B3)"

# Update selection to B3
$ws.Range("B3").Select()
